$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.667.78"
$ws.Range("E2").Value = "'  +5.47%  "
$ws.Range("D3").Value = "'3.068.38"
$ws.Range("E3").Value = "'  +3.62%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'551.46"
$ws.Range("E5").Value = "'  +6.05%  "
$ws.Range("D6").Value = "'139.41"
$ws.Range("E6").Value = "'  +8.02%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("D8").Value = "'3.060.35"
$ws.Range("E8").Value = "'  +3.44%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "'  +3.90%  "
$ws.Range("D10").Value = "'6.25"
$ws.Range("E10").Value = "'  +2.44%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "'  +2.82%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "'  +5.12%  "
$ws.Range("E13").Value = "'  +5.30%  "
$ws.Range("D14").Value = "'34.86"
$ws.Range("E14").Value = "'  +6.46%  "
$ws.Range("D15").Value = "'3.570.03"
$ws.Range("E15").Value = "'  +3.75%  "
$ws.Range("D16").Value = "'63.773.24"
$ws.Range("E16").Value = "'  +5.30%  "
$ws.Range("D17").Value = "'3.077.80"
$ws.Range("E17").Value = "'  +3.69%  "
$ws.Range("E18").Value = "'  -0.57%  "
$ws.Range("D19").Value = "'6.75"
$ws.Range("E19").Value = "'  +5.15%  "
$ws.Range("D20").Value = "'481.49"
$ws.Range("E20").Value = "'  +6.43%  "
$ws.Range("D21").Value = "'13.62"
$ws.Range("E21").Value = "'  +5.82%  "
$ws.Range("D22").Value = "'0.682"
$ws.Range("E22").Value = "'  +3.09%  "
$ws.Range("D23").Value = "'7.21"
$ws.Range("E23").Value = "'  +7.15%  "
$ws.Range("D24").Value = "'81.57"
$ws.Range("E24").Value = "'  +5.49%  "
$ws.Range("D25").Value = "'12.57"
$ws.Range("E25").Value = "'  +8.07%  "
$ws.Range("E26").Value = "'  +0.11%  "
$ws.Range("D27").Value = "'2.76"
$ws.Range("E27").Value = "'  +6.25%  "
$ws.Range("D28").Value = "'7.99"
$ws.Range("E28").Value = "'  +4.91%  "
$ws.Range("D29").Value = "'1.99"
$ws.Range("E29").Value = "'  +9.95%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "'  -0.22%  "
$ws.Range("D31").Value = "'26.06"
$ws.Range("E31").Value = "'  +4.80%  "
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "'  +1.95%  "
$ws.Range("D33").Value = "'2.44"
$ws.Range("E33").Value = "'  +9.60%  "
$ws.Range("D34").Value = "'5.70"
$ws.Range("E34").Value = "'  +8.08%  "
$ws.Range("D35").Value = "'55.83"
$ws.Range("E35").Value = "'  +1.52%  "
$ws.Range("D36").Value = "'6.01"
$ws.Range("E36").Value = "'  +5.93%  "
$ws.Range("D37").Value = "'469.24"
$ws.Range("E37").Value = "'  +5.07%  "
$ws.Range("D38").Value = "'0.0817"
$ws.Range("E38").Value = "'  +6.49%  "
$ws.Range("D39").Value = "'3.152.82"
$ws.Range("E39").Value = "'  +0.13%  "
$ws.Range("D40").Value = "'0.0397"
$ws.Range("E40").Value = "'  +6.37%  "
$ws.Range("E41").Value = "'  +4.24%  "
$ws.Range("D42").Value = "'8.27"
$ws.Range("E42").Value = "'  +4.49%  "
$ws.Range("D43").Value = "'2.62"
$ws.Range("E43").Value = "'  +9.86%  "
$ws.Range("D44").Value = "'28.17"
$ws.Range("E44").Value = "'  +12.69%  "
$ws.Range("D45").Value = "'0.253"
$ws.Range("E45").Value = "'  +5.43%  "
$ws.Range("E47").Value = "'  +8.19%  "
$ws.Range("D48").Value = "'0.109"
$ws.Range("E48").Value = "'  +2.82%  "
$ws.Range("D49").Value = "'0.0₃0514"
$ws.Range("E49").Value = "'  +2.90%  "
$ws.Range("E50").Value = "'  -0.55%  "
$ws.Range("D51").Value = "'2.07"
$ws.Range("E51").Value = "'  +7.76%  "
